# Realestate Update resale numbers 2023-06-08 14:44
# Appends a new data row (row 31) to the CityResaleNum sheet with the
# latest resale-number snapshot, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Text columns (Date / Time / Weekday / Week).
# A leading apostrophe forces these to be stored as text instead of being
# auto-parsed into date/number values, and ClearFormats() strips the
# "quote prefix" style Excel would otherwise tag the cell with, so the
# cells end up with plain default formatting - same as all the other rows.
$ws.Range("A$row").Value = "'2023-06-08"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = "'14:38:27"
$ws.Range("B$row").ClearFormats()

$ws.Range("C$row").Value = "'Thursday"
$ws.Range("C$row").ClearFormats()

$ws.Range("D$row").Value = "'23"
$ws.Range("D$row").ClearFormats()

# Numeric columns (city resale numbers).
$ws.Range("E$row").Value = 119058
$ws.Range("F$row").Value = 134346
$ws.Range("G$row").Value = 160254
$ws.Range("H$row").Value = 131292
$ws.Range("I$row").Value = 175565
$ws.Range("J$row").Value = 113097
$ws.Range("K$row").Value = 201045
$ws.Range("L$row").Value = 221028
$ws.Range("M$row").Value = 172789
$ws.Range("N$row").Value = 120077
$ws.Range("O$row").Value = 38607
$ws.Range("P$row").Value = 34472
$ws.Range("Q$row").Value = 50807
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36661
$ws.Range("T$row").Value = -1
